# Se mejora el menú
# Append new chat rows (34-46) to the "Chats" worksheet, mirroring the
# existing Fecha/Mensaje log layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the existing data (row 33 is the last one).
# Each entry is Fecha (date/time text), Mensaje (message text).
$rows = @(
    @("03-05-2022 09:39", "Hola"),
    @("03-05-2022 09:39", "hola"),
    @("03-05-2022 09:42", "hola"),
    @("03-05-2022 09:42", "1"),
    @("03-05-2022 09:42", "hola"),
    @("03-05-2022 09:42", "hola"),
    @("03-05-2022 09:42", "hola"),
    @("03-05-2022 09:48", "hola"),
    @("03-05-2022 09:48", "1"),
    @("03-05-2022 09:49", "1"),
    @("03-05-2022 09:51", "hola"),
    @("03-05-2022 09:52", "1"),
    @("03-05-2022 09:53", "hola")
)

$startRow = 34
$endRow = $startRow + $rows.Length - 1

# Pre-format column B of the new rows as Text so purely numeric-looking
# messages (e.g. "1") are stored as text rather than being converted to
# numbers by Excel's input parser.
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $fecha = $rows[$i][0]
    $mensaje = $rows[$i][1]
    $ws.Cells.Item($r, 1).Value = $fecha
    $ws.Cells.Item($r, 2).Value = $mensaje
}

# Restore the cells to the workbook's default (unformatted) style so the
# newly written rows look just like the pre-existing ones.
$ws.Range("B$startRow`:B$endRow").ClearFormats()
